# Rename CPUSIM into start script / re-arrange code into functions.
# The workbook gained a new simulation cycle ("cycle 11"), recorded as a
# new column M on both Sheet1 (per-address memory trace) and Sheet2
# (per-register summary). A couple of stale "0" placeholder cells in
# column B of Sheet1 (rows where a value changed but the label cell was
# never refreshed) are also brought back in sync with column C.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---- Sheet1 : add "cycle 11" as column M -------------------------------
$ws1.Cells.Item(1, 13).Value = "cycle 11"

for ($r = 2; $r -le 257; $r++) {
    $lVal = $ws1.Cells.Item($r, 12).Value2
    $ws1.Cells.Item($r, 13).Value = $lVal
}

# Row 4 (address 2) keeps its L column shown in hex ("0F") but the new
# cycle column continues the plain numeric series like the rest of the row.
$ws1.Cells.Item(4, 13).Value = 3

# ---- Sheet1 : fix stale column B labels so they track column C ---------
$ws1.Cells.Item(4, 2).Value = "3"
$ws1.Cells.Item(256, 2).Value = "5F"

# ---- Sheet2 : add "cycle 11" header as column M -------------------------
$ws2.Cells.Item(1, 13).Value = "cycle 11"
